$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.964.71'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.47%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.663.29'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -2.08%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.50'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -2.06%  '

$ws.Range("E6").Value = '  -1.28%  '

$ws.Range("E7").Value = '  +0.07%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2629'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.76%  '

$ws.Range("E9").Value = '  +1.76%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.81'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07403'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.83%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.669.61'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -1.60%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.486'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5803'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008427'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.20%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.96'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -2.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.022.79'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.40%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.899'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.80%  '

$ws.Range("E19").Value = '  +0.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.64'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -2.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '188.60'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.85%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.174'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -1.69%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '144.73'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.21%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.566'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.55%  '

$ws.Range("E26").Value = '  +3.14%  '

$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06584'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +15.80%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.308'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -1.35%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.312'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.92%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.519'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.00%  '

$ws.Range("E32").Value = '  -0.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.624'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -2.39%  '

$ws.Range("E34").Value = '  -1.08%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6045'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.29%  '

$ws.Range("E36").Value = '  -0.21%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.678'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -0.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.196'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +5.45%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01602'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.02%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.074.34'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -2.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8575'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.61%  '

$ws.Range("E42").Value = '  +0.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '100.28'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +1.19%  '

$ws.Range("B44").Value = 'RocketPoolETH'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.811.52'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.48%  '

$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000114'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +5.22%  '

$ws.Range("E46").Value = '  -1.40%  '

$ws.Range("E47").Value = '  -0.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.998'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05204'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4294'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -0.70%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.927'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.62%  '
